# PO Test Case Migration
#
# Update the QuantityBefore / QuantityAfter "expected value" cells on the
# InventoryRequisition test-data sheet with the migrated PO test-case
# values. The source data stores these numeric-looking IDs as literal text
# (e.g. "1559.0"), so each cell is forced to Text format before the write
# and then has that format cleared again afterwards -- this keeps the
# value stored as a shared string (matching the existing column contents)
# without leaving a stray per-cell number format behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InventoryRequisition")

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value() = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("H2") "1559.0"
Set-TextValue $ws.Range("I2") "1561.0"
Set-TextValue $ws.Range("H3") "910.0"
Set-TextValue $ws.Range("I3") "912.0"
